$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 2 ("童子杰") to make room
# for the two new people, shifting the existing data rows down.
$ws.Rows("2:3").Insert()

# The insert copies the header row's (bold/shaded) formatting down into
# the newly inserted rows; strip it back to the plain/default style used
# by the other data rows.
$ws.Rows("2:3").ClearFormats()

# New row 2: 俞鸿泰, 24
$ws.Range("A2").Value = "俞鸿泰"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "24"
$ws.Range("C2").Value = "牛马学院"
$ws.Range("D2").Value = "Connection refused: no further information: localhost/127.0.0.1:8081"

# New row 3: 七七, 23
$ws.Range("A3").Value = "七七"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "23"
$ws.Range("C3").Value = "牛马学院"
$ws.Range("D3").Value = "Connection refused: no further information: localhost/127.0.0.1:8081"
